$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added: modality description in table modalities - update the description of dataset 3
$ws.Range("G4").Value = "description of dataset 3, with speacial html l'ike > or & or < d'es fois"

# Update the active cell selection from N7 to I7
$ws.Range("I7").Select()
